# Add a LOT questions
#
# The "Drawdown" question (previously the last row of the top block, row
# 13) moves up to row 2 and gains two extra reference links (columns G
# and H). Every other question between the old row 2 and row 12 shifts
# down by one row to make room; everything from row 14 down is untouched.
#
# Cell-level hyperlinks in this workbook are anchored to a fixed address
# rather than following their cell's content, so instead of relying on a
# row insert/delete (which would leave stale hyperlink anchors behind),
# every affected cell's text is rewritten in place and then the whole
# hyperlink set for the sheet is rebuilt from scratch in its final,
# correct layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithm")

# ---------------------------------------------------------------------
# 0. Clear cells whose column is occupied "before" but not "after" for
#    a given row, since rewriting A/D/E in place would otherwise leave
#    stale leftovers from whatever used to live in that row.
# ---------------------------------------------------------------------
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()
$ws.Range("A5").Clear()
$ws.Range("E5").Clear()
$ws.Range("F11").Clear()
$ws.Range("F13").Clear()

# ---------------------------------------------------------------------
# 1. Rewrite the text/number content of rows 2-13 in place.
# ---------------------------------------------------------------------

# Row 2: now the "Drawdown" question (was row 13), with Freq bumped 2 -> 3
# and two brand new links appended in columns G and H.
$ws.Range("A2").Value2 = "Drawdown "
$ws.Range("D2").Value2 = 3
$ws.Range("E2").Value2 = "https://www.1point3acres.com/bbs/thread-523645-1-1.html"
$ws.Range("F2").Value2 = "https://www.1point3acres.com/bbs/thread-515407-1-1.html"
$ws.Range("G2").Value2 = "https://www.1point3acres.com/bbs/thread-545456-1-1.html"
$ws.Range("H2").Value2 = "https://leetcode.com/discuss/interview-question/125037/The-2-player-game-of-Drawdown-with-N-groups-of-stones"

# Row 3: the "finite state machine" question (was row 2).
$ws.Range("A3").Value2 = "设计一个有限状态机，如何给出一系列操作，有限状态机返回结果"
$ws.Range("D3").Value2 = 1
$ws.Range("E3").Value2 = "https://www.1point3acres.com/bbs/thread-544516-1-1.html"

# Row 4: (was row 3, unchanged content).
$ws.Range("A4").Value2 = "抓动物算积分游戏"
$ws.Range("D4").Value2 = 1
$ws.Range("E4").Value2 = "https://www.1point3acres.com/bbs/thread-538176-1-1.html"

# Row 5: (was row 4, unchanged content).
$ws.Range("B5").Value2 = 105
$ws.Range("C5").Value2 = "Construct Binary Tree from Preorder and Inorder Traversal   "
$ws.Range("D5").Value2 = 3

# Row 6: (was row 5, unchanged content).
$ws.Range("A6").Value2 = "给一个函数(只有+-，没有乘除)，让你写程序实现"
$ws.Range("D6").Value2 = 1
$ws.Range("E6").Value2 = "https://www.1point3acres.com/bbs/thread-536645-1-1.html"

# Row 7: (was row 6, unchanged content).
$ws.Range("A7").Value2 = "print histogram bar graph"
$ws.Range("D7").Value2 = 1
$ws.Range("E7").Value2 = "https://www.1point3acres.com/bbs/thread-535384-1-1.html"

# Row 8: (was row 7, unchanged content).
$ws.Range("A8").Value2 = "valid playlist"
$ws.Range("D8").Value2 = 2
$ws.Range("E8").Value2 = "https://www.1point3acres.com/bbs/thread-535184-1-1.html"

# Row 9: (was row 8, unchanged content).
$ws.Range("A9").Value2 = "connectFour game"
$ws.Range("D9").Value2 = 1
$ws.Range("E9").Value2 = "https://www.1point3acres.com/bbs/thread-528959-1-1.html"

# Row 10: (was row 9, unchanged content).
$ws.Range("A10").Value2 = "describeRelationship"
$ws.Range("D10").Value2 = 2
$ws.Range("E10").Value2 = "https://www.1point3acres.com/bbs/thread-527133-1-1.html"

# Row 11: (was row 10, unchanged content).
$ws.Range("A11").Value2 = "按7段数码管格式打印出数字"
$ws.Range("D11").Value2 = 1
$ws.Range("E11").Value2 = "https://www.1point3acres.com/bbs/thread-526112-1-1.html"

# Row 12: (was row 11, unchanged content).
$ws.Range("A12").Value2 = "cat query"
$ws.Range("D12").Value2 = 1
$ws.Range("E12").Value2 = "https://www.1point3acres.com/bbs/thread-525809-1-1.html"
$ws.Range("F12").Value2 = "https://leetcode.com/discuss/interview-question/124941/Square-or-Phone-screen-or-Search-cats"

# Row 13: (was row 12, unchanged content).
$ws.Range("A13").Value2 = "sales 分钱"
$ws.Range("D13").Value2 = 1
$ws.Range("E13").Value2 = "https://www.1point3acres.com/bbs/thread-524237-1-1.html"

# ---------------------------------------------------------------------
# 2. Rebuild every hyperlink on the sheet in its final, correct place.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @(
    @("E2",  "https://www.1point3acres.com/bbs/thread-523645-1-1.html"),
    @("F2",  "https://www.1point3acres.com/bbs/thread-515407-1-1.html"),
    @("G2",  "https://www.1point3acres.com/bbs/thread-545456-1-1.html"),
    @("H2",  "https://leetcode.com/discuss/interview-question/125037/The-2-player-game-of-Drawdown-with-N-groups-of-stones"),
    @("E3",  "https://www.1point3acres.com/bbs/thread-544516-1-1.html"),
    @("E4",  "https://www.1point3acres.com/bbs/thread-538176-1-1.html"),
    @("E6",  "https://www.1point3acres.com/bbs/thread-536645-1-1.html"),
    @("E7",  "https://www.1point3acres.com/bbs/thread-535384-1-1.html"),
    @("E8",  "https://www.1point3acres.com/bbs/thread-535184-1-1.html"),
    @("E9",  "https://www.1point3acres.com/bbs/thread-528959-1-1.html"),
    @("E10", "https://www.1point3acres.com/bbs/thread-527133-1-1.html"),
    @("E11", "https://www.1point3acres.com/bbs/thread-526112-1-1.html"),
    @("E12", "https://www.1point3acres.com/bbs/thread-525809-1-1.html"),
    @("F12", "https://leetcode.com/discuss/interview-question/124941/Square-or-Phone-screen-or-Search-cats"),
    @("E13", "https://www.1point3acres.com/bbs/thread-524237-1-1.html"),
    @("E14", "https://www.1point3acres.com/bbs/thread-521654-1-1.html"),
    @("E15", "https://www.1point3acres.com/bbs/thread-521588-1-1.html"),
    @("E16", "https://www.1point3acres.com/bbs/thread-521588-1-1.html"),
    @("E17", "https://www.1point3acres.com/bbs/thread-519586-1-1.html"),
    @("E18", "https://www.1point3acres.com/bbs/thread-519256-1-1.html"),
    @("E20", "https://www.1point3acres.com/bbs/thread-519256-1-1.html"),
    @("E21", "https://www.1point3acres.com/bbs/thread-515859-1-1.html"),
    @("E22", "https://www.1point3acres.com/bbs/thread-515859-1-1.html"),
    @("E23", "https://www.1point3acres.com/bbs/thread-476708-1-1.html"),
    @("E24", "https://www.1point3acres.com/bbs/thread-490256-1-1.html"),
    @("E25", "https://www.1point3acres.com/bbs/thread-490256-1-1.html"),
    @("E26", "https://www.1point3acres.com/bbs/thread-483851-1-1.html"),
    @("E27", "https://www.1point3acres.com/bbs/thread-481314-1-1.html"),
    @("E28", "https://www.1point3acres.com/bbs/thread-481314-1-1.html"),
    @("E29", "https://www.1point3acres.com/bbs/thread-479313-1-1.html"),
    @("E30", "https://www.1point3acres.com/bbs/thread-479313-1-1.html"),
    @("E31", "https://www.1point3acres.com/bbs/thread-479313-1-1.html"),
    @("E32", "https://www.1point3acres.com/bbs/thread-475381-1-1.html"),
    @("E33", "https://www.1point3acres.com/bbs/thread-475302-1-1.html"),
    @("E34", "https://www.1point3acres.com/bbs/thread-475302-1-1.html"),
    @("E35", "https://www.1point3acres.com/bbs/thread-474814-1-1.html")
)

foreach ($link in $links) {
    $ref = $link[0]
    $target = $link[1]
    $ws.Hyperlinks.Add($ws.Range($ref), $target) | Out-Null
    $ws.Range($ref).Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 3. Restore the recorded selections / active sheet from the diff.
# ---------------------------------------------------------------------
$ws.Range("J17").Select()

$ws2 = $wb.Worksheets.Item("design")
$ws2.Range("F11").Select()
$ws2.Activate()
